$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F7").Clear()
$ws.Range("H7").Clear()
